$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 37362
$ws.Range("E2").Value = 273
$ws.Range("F2").Value = 275
$ws.Range("G2").Value = 172
$ws.Range("H2").Value = 116
$ws.Range("I2").Value = 165
$ws.Range("J2").Value = -49
$ws.Range("K2").Value = 36598
$ws.Range("L2").Value = 23416
$ws.Range("M2").Value = 13182
$ws.Range("N2").Value = 11788
$ws.Range("O2").Value = 1394
$ws.Range("P2").Value = 203
$ws.Range("Q2").Value = -614
$ws.Range("R2").Value = -1877
$ws.Range("S2").Value = 1098
$ws.Range("T2").Value = 2285
$ws.Range("U2").Value = -2899
$ws.Range("V2").Value = 11834
$ws.Range("W2").Value = 0.73
$ws.Range("X2").Value = 0.31
$ws.Range("Y2").Value = 1.4
$ws.Range("Z2").Value = 0.36
$ws.Range("AA2").Value = 177.64
$ws.Range("AB2").Value = 5972.4
$ws.Range("AC2").Value = 4069
$ws.Range("AD2").Value = 34.16
$ws.Range("AE2").Value = 330810
$ws.Range("AF2").Value = 0.42
$ws.Range("AG2").Value = 2000
$ws.Range("AH2").Value = 1.44
$ws.Range("AI2").Value = 43.19
$ws.Range("AJ2").Value = 4055025
$ws.Range("D3").Value = 36679
$ws.Range("E3").Value = 891
$ws.Range("F3").Value = 891
$ws.Range("G3").Value = 408
$ws.Range("H3").Value = 372
$ws.Range("I3").Value = 371
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 35435
$ws.Range("L3").Value = 21945
$ws.Range("M3").Value = 13491
$ws.Range("N3").Value = 12104
$ws.Range("O3").Value = 1387
$ws.Range("P3").Value = 203
$ws.Range("Q3").Value = 1940
$ws.Range("R3").Value = -2666
$ws.Range("S3").Value = 403
$ws.Range("T3").Value = 1293
$ws.Range("U3").Value = 647
$ws.Range("V3").Value = 12347
$ws.Range("W3").Value = 2.43
$ws.Range("X3").Value = 1.01
$ws.Range("Y3").Value = 3.11
$ws.Range("Z3").Value = 1.03
$ws.Range("AA3").Value = 162.67
$ws.Range("AB3").Value = 6101.89
$ws.Range("AC3").Value = 9150
$ws.Range("AD3").Value = 11.2
$ws.Range("AE3").Value = 339663
$ws.Range("AF3").Value = 0.3
$ws.Range("AG3").Value = 2500
$ws.Range("AH3").Value = 2.44
$ws.Range("AI3").Value = 24.01
$ws.Range("AJ3").Value = 4055025
$ws.Range("D4").Value = 30601
$ws.Range("E4").Value = 615
$ws.Range("F4").Value = 606
$ws.Range("G4").Value = 385
$ws.Range("H4").Value = 201
$ws.Range("I4").Value = 345
$ws.Range("J4").Value = -145
$ws.Range("K4").Value = 36902
$ws.Range("L4").Value = 23272
$ws.Range("M4").Value = 13630
$ws.Range("N4").Value = 12388
$ws.Range("O4").Value = 1242
$ws.Range("P4").Value = 203
$ws.Range("Q4").Value = 3299
$ws.Range("R4").Value = -1944
$ws.Range("S4").Value = -651
$ws.Range("T4").Value = 679
$ws.Range("U4").Value = 2620
$ws.Range("V4").Value = 11976
$ws.Range("W4").Value = 2.01
$ws.Range("X4").Value = 0.66
$ws.Range("Y4").Value = 2.82
$ws.Range("Z4").Value = 0.55
$ws.Range("AA4").Value = 170.74
$ws.Range("AB4").Value = 6241.31
$ws.Range("AC4").Value = 8519
$ws.Range("AD4").Value = 11.92
$ws.Range("AE4").Value = 347648
$ws.Range("AF4").Value = 0.29
$ws.Range("AG4").Value = 2500
$ws.Range("AH4").Value = 2.46
$ws.Range("AI4").Value = 25.79
$ws.Range("AJ4").Value = 4055025
$ws.Range("D5").Value = 32951
$ws.Range("E5").Value = 611
$ws.Range("F5").Value = 611
$ws.Range("G5").Value = 221
$ws.Range("H5").Value = 51
$ws.Range("I5").Value = 132
$ws.Range("J5").Value = -81
$ws.Range("K5").Value = 37043
$ws.Range("L5").Value = 23484
$ws.Range("M5").Value = 13558
$ws.Range("N5").Value = 12396
$ws.Range("O5").Value = 1162
$ws.Range("P5").Value = 203
$ws.Range("Q5").Value = 1001
$ws.Range("R5").Value = -107
$ws.Range("S5").Value = -460
$ws.Range("T5").Value = 1343
$ws.Range("U5").Value = -341
$ws.Range("V5").Value = 11468
$ws.Range("W5").Value = 1.86
$ws.Range("X5").Value = 0.16
$ws.Range("Y5").Value = 1.07
$ws.Range("Z5").Value = 0.14
$ws.Range("AA5").Value = 173.21
$ws.Range("AB5").Value = 6257.78
$ws.Range("AC5").Value = 3259
$ws.Range("AD5").Value = 34.67
$ws.Range("AE5").Value = 347883
$ws.Range("AF5").Value = 0.32
$ws.Range("AG5").Value = 3000
$ws.Range("AH5").Value = 2.65
$ws.Range("AI5").Value = 80.89
$ws.Range("AJ5").Value = 4055025
$ws.Range("D6").Value = 34581
$ws.Range("E6").Value = 823
$ws.Range("F6").Value = 823
$ws.Range("G6").Value = 461
$ws.Range("H6").Value = 266
$ws.Range("I6").Value = 321
$ws.Range("K6").Value = 36994
$ws.Range("L6").Value = 23311
$ws.Range("M6").Value = 13683
$ws.Range("N6").Value = 12577
$ws.Range("P6").Value = 203
$ws.Range("Q6").Value = 1250
$ws.Range("R6").Value = -1214
$ws.Range("S6").Value = -340
$ws.Range("T6").Value = 1415
$ws.Range("U6").Value = -165
$ws.Range("V6").Value = 11161
$ws.Range("W6").Value = 2.38
$ws.Range("X6").Value = 0.77
$ws.Range("Y6").Value = 2.57
$ws.Range("Z6").Value = 0.72
$ws.Range("AA6").Value = 170.36
$ws.Range("AB6").Value = 6338.81
$ws.Range("AC6").Value = 7905
$ws.Range("AD6").Value = 11.51
$ws.Range("AE6").Value = 352950
$ws.Range("AF6").Value = 0.26
$ws.Range("AG6").Value = 3000
$ws.Range("AH6").Value = 3.3
$ws.Range("AI6").Value = 33.35
$ws.Range("AJ6").Value = 4055025
$ws.Range("D7").Value = 35703
$ws.Range("E7").Value = 910
$ws.Range("G7").Value = 714
$ws.Range("H7").Value = 546
$ws.Range("I7").Value = 556
$ws.Range("K7").Value = 37278
$ws.Range("L7").Value = 23146
$ws.Range("M7").Value = 14133
$ws.Range("N7").Value = 13029
$ws.Range("P7").Value = 203
$ws.Range("Q7").Value = 863
$ws.Range("R7").Value = -1830
$ws.Range("S7").Value = 517
$ws.Range("T7").Value = 1192
$ws.Range("U7").Value = 522
$ws.Range("W7").Value = 2.55
$ws.Range("X7").Value = 1.53
$ws.Range("Y7").Value = 4.35
$ws.Range("Z7").Value = 1.47
$ws.Range("AA7").Value = 163.77
$ws.Range("AC7").Value = 13720
$ws.Range("AD7").Value = 5.72
$ws.Range("AE7").Value = 365633
$ws.Range("AF7").Value = 0.21
$ws.Range("AG7").Value = 3000
$ws.Range("AH7").Value = 3.82
$ws.Range("AI7").Value = 21.87
$ws.Range("D8").Value = 36526
$ws.Range("E8").Value = 890
$ws.Range("G8").Value = 682
$ws.Range("H8").Value = 522
$ws.Range("I8").Value = 547
$ws.Range("K8").Value = 37668
$ws.Range("L8").Value = 23128
$ws.Range("M8").Value = 14540
$ws.Range("N8").Value = 13442
$ws.Range("P8").Value = 203
$ws.Range("Q8").Value = 1678
$ws.Range("R8").Value = -1264
$ws.Range("S8").Value = 52
$ws.Range("T8").Value = 1120
$ws.Range("U8").Value = 420
$ws.Range("W8").Value = 2.44
$ws.Range("X8").Value = 1.43
$ws.Range("Y8").Value = 4.13
$ws.Range("Z8").Value = 1.39
$ws.Range("AA8").Value = 159.06
$ws.Range("AC8").Value = 13489
$ws.Range("AD8").Value = 5.82
$ws.Range("AE8").Value = 377223
$ws.Range("AF8").Value = 0.21
$ws.Range("AG8").Value = 3000
$ws.Range("AH8").Value = 3.82
$ws.Range("AI8").Value = 22.24
$ws.Range("D9").Value = 38215
$ws.Range("E9").Value = 903
$ws.Range("G9").Value = 711
$ws.Range("H9").Value = 543
$ws.Range("I9").Value = 567
$ws.Range("K9").Value = 38438
$ws.Range("L9").Value = 23468
$ws.Range("M9").Value = 14970
$ws.Range("N9").Value = 13923
$ws.Range("P9").Value = 203
$ws.Range("Q9").Value = 1761
$ws.Range("R9").Value = -1026
$ws.Range("S9").Value = -172
$ws.Range("T9").Value = 1120
$ws.Range("U9").Value = 741
$ws.Range("W9").Value = 2.36
$ws.Range("X9").Value = 1.42
$ws.Range("Y9").Value = 4.14
$ws.Range("Z9").Value = 1.43
$ws.Range("AA9").Value = 156.77
$ws.Range("AC9").Value = 13983
$ws.Range("AD9").Value = 5.61
$ws.Range("AE9").Value = 390722
$ws.Range("AF9").Value = 0.2
$ws.Range("AG9").Value = 3000
$ws.Range("AH9").Value = 3.82
$ws.Range("AI9").Value = 21.45
